$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RES_RAPIDA (column I) and TOTAL_PR (column J) values
# to reflect the "trained personeel in the last 2 years" question.

$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 57

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 38

$ws.Range("I4").Value = 12
$ws.Range("J4").Value = 34

$ws.Range("I5").Value = 12
$ws.Range("J5").Value = 37

$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 41

$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 30
